$d = $word.ActiveDocument

# Edit 1: CI pipeline sentence
$d.Content.Find.Execute(
    "Improved the product’s delivery speed by increasing the CI pipeline’s success rate by 10%",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Improved the product’s delivery speed by 10% by fixing the non-deterministic errors in the CI pipeline", 2)

# Edit 2: low latency trading sentence
$d.Content.Find.Execute(
    "low latency and blocking others’ transactions in EVM-compatible",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "low-latency trading in EVM-compatible", 2)
